$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fix B31 value (timestamp re-computed on re-run)
$ws.Range("B31").Value = 43417.93120255787

# Copy row 31 (template formatting/styles) down into the 7 new rows (32-38)
$srcRow = $ws.Range("A31:N31")
for ($r = 32; $r -le 38; $r++) {
    $dstRow = $ws.Range("A$r" + ":N$r")
    $srcRow.Copy($dstRow)
}

# Row 32
$ws.Range("A32").Value = 30
$ws.Range("B32").Value = 43417.93338178241
$ws.Range("C32").Value = "JAFFE"
$ws.Range("D32").Value = "LBP8"
$ws.Range("E32").Value = "LDA1"
$ws.Range("F32").Value = "SVM2"
$ws.Range("G32").Value = 0.4553990610328639
$ws.Range("H32").Value = 18.10745525360107
$ws.Range("I32").Value = 9.070727109909058
$ws.Range("J32").Value = 8.942731380462646
$ws.Range("K32").Value = 0.07999753952026367
$ws.Range("L32").Value = 0.005600118637084961
$ws.Range("M32").Value = 0.001500535011291504
$ws.Range("N32").Value = 0.0003999948501586914

# Row 33
$ws.Range("A33").Value = 31
$ws.Range("B33").Value = 43417.93377541667
$ws.Range("C33").Value = "JAFFE"
$ws.Range("D33").Value = "LBP8"
$ws.Range("E33").Value = "PCA1"
$ws.Range("F33").Value = "SVM2"
$ws.Range("G33").Value = 0.4647887323943662
$ws.Range("H33").Value = 18.38443160057068
$ws.Range("I33").Value = 9.267715930938721
$ws.Range("J33").Value = 8.999729156494141
$ws.Range("K33").Value = 0.1039659976959229
$ws.Range("L33").Value = 0.002099084854125977
$ws.Range("M33").Value = 0.007399535179138184
$ws.Range("N33").Value = 0.0005000114440917968

# Row 34
$ws.Range("A34").Value = 32
$ws.Range("B34").Value = 43417.93413646991
$ws.Range("C34").Value = "JAFFE"
$ws.Range("D34").Value = "LBP8"
$ws.Range("E34").Value = "LDA1"
$ws.Range("F34").Value = "SVM1"
$ws.Range("G34").Value = 0.5070422535211268
$ws.Range("H34").Value = 18.65042853355408
$ws.Range("I34").Value = 9.126715183258057
$ws.Range("J34").Value = 9.463715076446533
$ws.Range("K34").Value = 0.04698753356933594
$ws.Range("L34").Value = 0.003199863433837891
$ws.Range("M34").Value = 0.001100015640258789
$ws.Range("N34").Value = 0.000099945068359375

# Row 35
$ws.Range("A35").Value = 33
$ws.Range("B35").Value = 43417.93451767361
$ws.Range("C35").Value = "JAFFE"
$ws.Range("D35").Value = "LBP8"
$ws.Range("E35").Value = "LDA1"
$ws.Range("F35").Value = "SVM3"
$ws.Range("G35").Value = 0.3051643192488263
$ws.Range("H35").Value = 18.46443343162537
$ws.Range("I35").Value = 9.061725616455078
$ws.Range("J35").Value = 9.333709955215454
$ws.Range("K35").Value = 0.05399847030639648
$ws.Range("L35").Value = 0.003000235557556152
$ws.Range("M35").Value = 0.001799774169921875
$ws.Range("N35").Value = 0.0001999139785766602

# Row 36
$ws.Range("A36").Value = 34
$ws.Range("B36").Value = 43417.93495306713
$ws.Range("C36").Value = "JAFFE"
$ws.Range("D36").Value = "LBP8"
$ws.Range("E36").Value = "PCA1"
$ws.Range("F36").Value = "SVM3"
$ws.Range("G36").Value = 0.5539906103286385
$ws.Range("H36").Value = 18.37643694877625
$ws.Range("I36").Value = 9.139714479446411
$ws.Range("J36").Value = 9.133739709854126
$ws.Range("K36").Value = 0.08898758888244629
$ws.Range("L36").Value = 0.001899242401123047
$ws.Range("M36").Value = 0.006199479103088379
$ws.Range("N36").Value = 0.0005001068115234375

# Row 37
$ws.Range("A37").Value = 35
$ws.Range("B37").Value = 43417.93529873843
$ws.Range("C37").Value = "JAFFE"
$ws.Range("D37").Value = "LBP8"
$ws.Range("E37").Value = "PCA1"
$ws.Range("F37").Value = "SVM1"
$ws.Range("G37").Value = 0.4694835680751174
$ws.Range("H37").Value = 18.76842522621155
$ws.Range("I37").Value = 9.079715728759766
$ws.Range("J37").Value = 9.563712596893311
$ws.Range("K37").Value = 0.1119873523712158
$ws.Range("L37").Value = 0.002299070358276367
$ws.Range("M37").Value = 0.008099508285522462
$ws.Range("N37").Value = 0.0004001379013061523

# Row 38
$ws.Range("A38").Value = 36
$ws.Range("B38").Value = 43417.93577364587
$ws.Range("C38").Value = "JAFFE"
$ws.Range("D38").Value = "LBP8"
$ws.Range("E38").Value = "PCA1"
$ws.Range("F38").Value = "SVM2"
$ws.Range("G38").Value = 0.4647887323943662
$ws.Range("H38").Value = 18.54843211174011
$ws.Range("I38").Value = 9.096726894378662
$ws.Range("J38").Value = 9.32970929145813
$ws.Range("K38").Value = 0.1079964637756348
$ws.Range("L38").Value = 0.001898932456970215
$ws.Range("M38").Value = 0.007799482345581055
$ws.Range("N38").Value = 0.0006002426147460937

Write-Output "done"